$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header "TOTAL" in H1, matching style of existing header (bold, border)
$ws.Range("A1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "TOTAL"

# New column width for H (closest representable value to the authored 10.42578125)
$ws.Columns.Item(8).ColumnWidth = 9.59

# New row 22: TOTAL label + sums, matching style of existing header (bold, border)
$ws.Range("A1").Copy()
$ws.Range("A22").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A22").Value = "TOTAL"

$ws.Range("F22").Formula = "=SUM(F2:F20)"
$ws.Range("G22").Formula = "=SUM(G2:G20)"
$ws.Range("H22").Formula = "=SUM(F22:G22)"

$excel.CutCopyMode = $false

# Move selection similar to author's final cursor position
[void]$ws.Range("F30").Select()

$wb.Save()
